# Added New Mac-Address and Document Types
# Appends 5 new "eng" rows (regcntr_id 10002 / machine_id 10032) for
# device_ids 3000176-3000180 to the bottom of the master data sheet,
# switches the workbook to manual calculation, and leaves the
# selection on the first new data row (column D) as the author did.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$regcntrId  = 10002
$machineId  = 10032
$firstDevId = 3000176
$rowCount   = 5
$startRow   = 157

for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $startRow + $i

    $ws.Cells.Item($r, 1).Value = $regcntrId
    $ws.Cells.Item($r, 2).Value = $machineId
    $ws.Cells.Item($r, 3).Value = $firstDevId + $i
    $ws.Cells.Item($r, 4).Value = "eng"
    $ws.Cells.Item($r, 5).Value = $true
    $ws.Cells.Item($r, 6).Value = "superadmin"
    $ws.Cells.Item($r, 7).Value = "now()"
    $ws.Cells.Item($r, 8).Value = "now()"
}

# Workbook switches from automatic to manual calculation.
$excel.Calculation = -4135

# Leave the active selection on the first cell of the newly entered data.
[void]$ws.Range("D157").Select()
